# add group min/max weights
#
# The sheet holds a ticker/weight table in columns A (Ticker) and B (Weight),
# rows 2..32 (row 1 is the header). This edit:
#  - re-baselines the whole ticker list/order (TLT and GLTR move to the top,
#    right under the header; EMB is dropped entirely; the remaining tickers
#    keep the same relative order but shift to fill the gaps),
#  - rewrites every group's Weight to its new min/max-derived target value,
#  - removes the now-unused last row (old row 32 / EMB+TLT leftover).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Ticker order for rows 2..31, and the new Weight for each of those rows.
$tickers = @(
    "TLT", "GLTR", "FDX", "WMT", "AMZN", "CORE", "BLK", "LMT", "ORCL", "NTRS",
    "TSM", "SJM", "MDLZ", "REGI", "V", "MSFT", "JNJ", "TPH", "VIRT", "AXP",
    "BX", "CNC", "LDOS", "MDT", "MRK", "NKE", "PFE", "SYY", "NRZ", "OHI"
)
$weights = @(
    0.1212, 0.15, 0.0304, 0.05, 0.1, 0.005, 0.05, 0.0232, 0.005, 0.005,
    0.05, 0.005, 0.005, 0.05, 0.05, 0.05, 0.005, 0.005, 0.005, 0.005,
    0.05, 0.0404, 0.05, 0.005, 0.0139, 0.05, 0.005, 0.0058, 0.005, 0.005
)

for ($i = 0; $i -lt $tickers.Length; $i++) {
    $row = $i + 2
    if ($ws.Cells.Item($row, 1).Value2 -ne $tickers[$i]) {
        $ws.Cells.Item($row, 1).Value = $tickers[$i]
    }
    if ($ws.Cells.Item($row, 2).Value2 -ne $weights[$i]) {
        $ws.Cells.Item($row, 2).Value = $weights[$i]
    }
}

# Drop the now-superfluous last row (31 data rows remain instead of 31+EMB=32).
$ws.Rows.Item(32).Delete()
